# Generate Report for Handback
# Updates the "Correspond Handoff Datetime" (E2) and
# "Correspond Handback DateTime" (H2) timestamps on the zh-cn and de-de
# sheets to reflect the new report generation time.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-20 17:00:14"
$wsZhCn.Range("H2").Value = "2016-03-20 17:00:46"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-20 17:00:18"
$wsDeDe.Range("H2").Value = "2016-03-20 17:00:51"
